$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.939.34"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "3.095.15"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'542.05"
$ws.Range("E5").Value = "  -2.03%  "

$ws.Range("D6").Value = "'137.00"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.088.26"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -6.60%  "

$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("E13").Value = "  +4.78%  "

$ws.Range("D14").Value = "'34.84"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").Value = "3.595.97"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "63.915.42"
$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D18").Value = "3.094.51"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").Value = "'490.67"
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'80.02"
$ws.Range("E24").Value = "  +3.73%  "

$ws.Range("D25").Value = "'12.32"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("D33").Value = "'2.43"
$ws.Range("E33").Value = "  -3.23%  "

$ws.Range("D34").Value = "'57.46"
$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("E35").Value = "  +5.09%  "

$ws.Range("D36").Value = "'497.77"
$ws.Range("E36").Value = "  -5.69%  "

$ws.Range("D37").Value = "'6.08"
$ws.Range("E37").Value = "  +3.44%  "

$ws.Range("D38").Value = "3.212.47"
$ws.Range("E38").Value = "  +5.72%  "

$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("D41").Value = "'0.119"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("E42").Value = "  +3.29%  "

$ws.Range("D43").Value = "'8.19"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("D44").Value = "'0.259"
$ws.Range("E44").Value = "  +2.59%  "

$ws.Range("E46").Value = "  +8.81%  "

$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "  +1.47%  "

$ws.Range("D48").Value = "'121.60"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").Value = "'24.85"
$ws.Range("E49").Value = "  +4.88%  "

$ws.Range("E50").Value = "  +3.39%  "

$ws.Range("E51").Value = "  -2.55%  "
